$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 18. This pushes the old row 18 (period 2505,
#    value 35338) down to row 19, and the trailing signature block (old rows
#    23-24) down to rows 24-25, while rows 16-17 remain untouched.
$ws.Rows("18:18").Insert()

# 2. Give the freshly inserted row 18 the same formatting as row 17 (the row
#    immediately above it) so it matches the rest of the detail table.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# 3. Populate the new row 18 with the data that used to be on row 16 (period
#    2507) -- this becomes a regular (non-bold) detail row now.
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1002058702"
$ws.Range("D18").Value = "MARIA JOSE PADILLA HERNANDEZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 88345
$ws.Range("G18").Value = 2208640

# 4. Row 16 now becomes period 2505, taking on the smaller overdue value that
#    used to sit on the old last row (18).
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 35338

# Row 17 (period 2506) is unchanged.

# 5. Row 19 (old row 18, shifted down by the insert) becomes the new final
#    period, 2508, with the same amounts as the other regular rows.
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 88345

# 6. Update the account summary figures: total overdue value and period
#    count (now 4 periods instead of 3).
$ws.Range("E11").Value = 300373
$ws.Range("F13").Value = 4
